$wb = $excel.ActiveWorkbook

# Generate Report for Handback
# For the zh-cn sheet, row 3 (the "ed31471e..." handoff row) gets its own
# distinct Handoff/Handback datetimes instead of sharing row 2's values.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-21 20:47:04"
$wsZhCn.Range("H3").Value = "2016-03-21 20:47:28"

# For the de-de sheet, row 3 (the "ed31471e..." handoff row) gets its own
# distinct Handoff/Handback datetimes instead of sharing row 2's values.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-21 20:47:08"
$wsDeDe.Range("H3").Value = "2016-03-21 20:47:34"
